$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D").Insert()

# Copy formatting from the (now-shifted) neighboring column E into the new column D
# for each of the three data blocks on this sheet (the gaps at rows 36/78 and the
# label-only rows 5,6,37,79 have no D:K cells and must be left untouched).
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new column D with the newest reporting periods figures
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 108500
$ws.Range("D9").Value = 22600
$ws.Range("D10").Value = 85900
$ws.Range("D12").Value = 19300
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 133500
$ws.Range("D18").Value = -25000
$ws.Range("D20").Value = 2100
$ws.Range("D21").Value = -21000
$ws.Range("D22").Value = "NA"
$ws.Range("D23").Value = -22900
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -22900
$ws.Range("D27").Value = -22900
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -2100
$ws.Range("D33").Value = -22900
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -22900
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 9500
$ws.Range("D42").Value = 91300
$ws.Range("D43").Value = 19600
$ws.Range("D44").Value = 11600
$ws.Range("D45").Value = 2700
$ws.Range("D46").Value = 134700
$ws.Range("D47").Value = "NA"
$ws.Range("D48").Value = 5900
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 400
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 141000
$ws.Range("D57").Value = 6200
$ws.Range("D58").Value = "NA"
$ws.Range("D59").Value = 13500
$ws.Range("D60").Value = 19700
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 200
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 20000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -187800
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 121000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -22900
$ws.Range("D83").Value = 1900
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -13800
$ws.Range("D91").Value = -2100
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -10000
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 13500
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -10400
